$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '''35.330.28'
$ws.Range("E2").Value = '  -0.92%  '
# Row 3
$ws.Range("D3").Value = '''1.904.14'
$ws.Range("E3").Value = '  -0.08%  '
# Row 4
$ws.Range("E4").Value = '  -0.13%  '
# Row 5
$ws.Range("D5").Value = '''0.692'
$ws.Range("E5").Value = '  +9.16%  '
# Row 6
$ws.Range("D6").Value = '''246.05'
$ws.Range("E6").Value = '  +0.40%  '
# Row 7
$ws.Range("E7").Value = '  -0.30%  '
# Row 8
$ws.Range("D8").Value = '''41.35'
$ws.Range("E8").Value = '  -3.39%  '
# Row 9
$ws.Range("E9").Value = '  +3.74%  '
# Row 10
$ws.Range("D10").Value = '''53.13'
$ws.Range("E10").Value = '  +11.63%  '
# Row 11
$ws.Range("D11").Value = '''0.0725'
$ws.Range("E11").Value = '  +2.38%  '
# Row 12
$ws.Range("E12").Value = '  +0.12%  '
# Row 13
$ws.Range("D13").Value = '''2.180.80'
$ws.Range("E13").Value = '  +0.08%  '
# Row 14
$ws.Range("D14").Value = '''12.40'
$ws.Range("E14").Value = '  -0.89%  '
# Row 15
$ws.Range("D15").Value = '''0.709'
$ws.Range("E15").Value = '  +2.45%  '
# Row 16
$ws.Range("D16").Value = '''1.905.91'
$ws.Range("E16").Value = '  +0.15%  '
# Row 17
$ws.Range("D17").Value = '''4.84'
$ws.Range("E17").Value = '  +0.66%  '
# Row 18
$ws.Range("D18").Value = '''35.346.87'
$ws.Range("E18").Value = '  -0.85%  '
# Row 19
$ws.Range("D19").Value = '''72.32'
$ws.Range("E19").Value = '  -0.03%  '
# Row 20
$ws.Range("D20").Value = '''0.0₃0821'
$ws.Range("E20").Value = '  +1.14%  '
# Row 21
$ws.Range("D21").Value = '''241.26'
$ws.Range("E21").Value = '  -1.47%  '
# Row 22
$ws.Range("D22").Value = '''12.72'
$ws.Range("E22").Value = '  +2.01%  '
# Row 23
$ws.Range("D23").Value = '''4.81'
$ws.Range("E23").Value = '  -1.81%  '
# Row 25
$ws.Range("E25").Value = '  +1.43%  '
# Row 26
$ws.Range("D26").Value = '''2.30'
$ws.Range("E26").Value = '  +9.60%  '
# Row 27
$ws.Range("D27").Value = '''168.44'
$ws.Range("E27").Value = '  -1.70%  '
# Row 28
$ws.Range("D28").Value = '''8.60'
$ws.Range("E28").Value = '  +2.22%  '
# Row 29
$ws.Range("E29").Value = '  +4.30%  '
# Row 30
$ws.Range("D30").Value = '''18.40'
$ws.Range("E30").Value = '  +2.10%  '
# Row 31
$ws.Range("D31").Value = '''4.138.76'
$ws.Range("E31").Value = '  +21.22%  '
# Row 32
$ws.Range("D32").Value = '''4.16'
$ws.Range("E32").Value = '  +1.41%  '
# Row 33
$ws.Range("B33").Value = 'Hedera'
$ws.Range("C33").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D33").Value = '''0.0571'
$ws.Range("E33").Value = '  +0.81%  '
# Row 34
$ws.Range("B34").Value = 'ImmutableX'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D34").Value = '''0.941'
$ws.Range("E34").Value = '  -1.59%  '
# Row 35
$ws.Range("E35").Value = '  -0.17%  '
# Row 36
$ws.Range("D36").Value = '''4.11'
$ws.Range("E36").Value = '  +0.14%  '
# Row 37
$ws.Range("D37").Value = '''1.79'
$ws.Range("E37").Value = '  +1.11%  '
# Row 38
$ws.Range("E38").Value = '  +8.17%  '
# Row 39
$ws.Range("D39").Value = '''2.02'
$ws.Range("E39").Value = '  -1.43%  '
# Row 40
$ws.Range("D40").Value = '''0.0666'
$ws.Range("E40").Value = '  +11.61%  '
# Row 41
$ws.Range("E41").Value = '  -0.96%  '
# Row 42
$ws.Range("E42").Value = '  +2.09%  '
# Row 43
$ws.Range("D43").Value = '''16.08'
$ws.Range("E43").Value = '  +4.41%  '
# Row 44
$ws.Range("D44").Value = '''89.93'
$ws.Range("E44").Value = '  -1.54%  '
# Row 45
$ws.Range("D45").Value = '''1.347.29'
$ws.Range("E45").Value = '  -1.58%  '
# Row 46
$ws.Range("D46").Value = '''2.44'
$ws.Range("E46").Value = '  +3.38%  '
# Row 47
$ws.Range("D47").Value = '''12.72'
$ws.Range("E47").Value = '  -3.28%  '
# Row 48
$ws.Range("B48").Value = 'HuobiToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D48").Value = '''2.42'
$ws.Range("E48").Value = '  +0.03%  '
# Row 49
$ws.Range("B49").Value = 'MXToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D49").Value = '''2.79'
$ws.Range("E49").Value = '  +0.78%  '
# Row 50
$ws.Range("B50").Value = 'MultiversX'
$ws.Range("C50").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D50").Value = '''45.75'
$ws.Range("E50").Value = '  -0.29%  '
# Row 51
$ws.Range("D51").Value = '''6.54'
$ws.Range("E51").Value = '  -3.12%  '
